$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3F")
$ws.Activate()

# Row 7, column C (Admission No) was recorded as text "05" from a form
# submission glitch; the sync normalizes it to the numeric value 5, like
# its siblings in the same column.
$ws.Cells.Item(7, 3).Value = 5

# New form submission synced in as row 8.
$ws.Cells.Item(8, 1).Value = "2026-02-08 22:16:43"
$ws.Cells.Item(8, 2).Value = "Usman Muhammad Gubio"
$ws.Cells.Item(8, 3).NumberFormat = "@"
$ws.Cells.Item(8, 3).Value = "05"
$ws.Cells.Item(8, 3).Style = "Normal"
$ws.Cells.Item(8, 4).Value = 9
